# edit.ps1 - Apply the recorded changes to aa_writing_02-style workbook
#
# Summary of edits performed by the original author:
#  1. Cells A53:A76 on Sheet1 were retyped to "ThinkPad L560" (matching A52),
#     which orphans the previously-unique "ThinkPad L561".."ThinkPad L584"
#     shared strings (Excel compacts/removes unused shared strings on save).
#  2. New raw benchmark numbers were typed into L77:N80 and F82:H85 on Sheet1.
#  3. The J77 formula was retyped directly as "=I77/D77" (fixing a stale
#     shared formula that had incorrectly been using AVERAGE(F:H)); the same
#     corrected formula was then carried down through J78:J86.
#  4. A new blank worksheet ("Sheet2") was added after Sheet1.
#  5. Selection / active-cell state was updated on both sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1. Retype A53:A76 to match A52 ("ThinkPad L560"); this collapses the
#    now-unused ThinkPad L561..L584 shared-string entries automatically
#    when the workbook is saved.
# ---------------------------------------------------------------------
$modelName = $ws1.Range("A52").Value2
$ws1.Range("A53:A76").Value = $modelName

# ---------------------------------------------------------------------
# 2. Fill in newly-collected benchmark figures.
# ---------------------------------------------------------------------
$ws1.Range("L77").Value = 1.427
$ws1.Range("M77").Value = 1.2470000000000001
$ws1.Range("N77").Value = 1.298

$ws1.Range("L78").Value = 144.94300000000001
$ws1.Range("M78").Value = 142.78899999999999
$ws1.Range("N78").Value = 144.928

$ws1.Range("L79").Value = 12871
$ws1.Range("M79").Value = 12825.7
$ws1.Range("N79").Value = 12847.1

$ws1.Range("L80").Value = 1284350
$ws1.Range("M80").Value = 1282130

$ws1.Range("F82").Value = 5.1123799999999999
$ws1.Range("G82").Value = 5.1094920000000004
$ws1.Range("H82").Value = 5.1664760000000003

$ws1.Range("F83").Value = 95.103114000000005
$ws1.Range("G83").Value = 103.511651
$ws1.Range("H83").Value = 95.155928000000003

$ws1.Range("F84").Value = 9171.7396669999998
$ws1.Range("G84").Value = 9165.2386989999995
$ws1.Range("H84").Value = 9157.5771719999993

$ws1.Range("F85").Value = 917663.71276000002

# ---------------------------------------------------------------------
# 3. Correct the J77 formula (was a mis-shared AVERAGE formula) and fill
#    the fix down through J86.
# ---------------------------------------------------------------------
$ws1.Range("J77").Formula = "=I77/D77"
for ($r = 78; $r -le 86; $r++) {
    $ws1.Range("J$r").Formula = "=I$r/D$r"
}

# ---------------------------------------------------------------------
# 4. Add the new blank "Sheet2" right after Sheet1.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "Sheet2"
$ws2.Activate()
$ws2.Range("H17").Select()

# ---------------------------------------------------------------------
# 5. Restore Sheet1 as the active sheet and update its selection state.
# ---------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("J85").Select()
